# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted into the "Piña" (pineapple)
# price table. It becomes the new row 209, pushing every existing
# record that used to live at row 209 (and all the rows after it)
# down by one. The sheet's used range grows from A1:T288 to A1:T289.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 209 - this shifts rows 209..288 down to 210..289
# and carries the date-number-format of column D forward automatically.
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A209").Value = 4
$ws.Range("B209").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C209").Value = "Los Lagos"
$ws.Range("D209").Value = 44795
$ws.Range("E209").Value = 10
$ws.Range("F209").Value = "Fruta"
$ws.Range("G209").Value = 100108
$ws.Range("H209").Value = "Tropicales y subtropicales"
$ws.Range("I209").Value = 100108005
$ws.Range("J209").Value = "Piña"
$ws.Range("K209").Value = "Caramelo"
$ws.Range("L209").Value = "Primera"
$ws.Range("M209").Value = 60
$ws.Range("N209").Value = 23000
$ws.Range("O209").Value = 23000
$ws.Range("P209").Value = 23000
$ws.Range("Q209").Value = "$/caja 14 unidades"
$ws.Range("R209").Value = "Ecuador"
$ws.Range("S209").Value = 1643
$ws.Range("T209").Value = 14
